$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Column 3 = "Name" (surname) -> widened by ~2mm: 1418 dxa (70.9pt) -> 1548 dxa (77.4pt)
$tbl.Columns.Item(3).Width = 77.4

# Column 4 = "Vorname" (first name) -> narrowed by ~2mm: 1418 dxa (70.9pt) -> 1288 dxa (64.4pt)
$tbl.Columns.Item(4).Width = 64.4

# The table's overall preferred width now reflects the fixed sum of the grid
# columns (11351 dxa = 567.55pt) instead of "auto".
$tbl.PreferredWidthType = 3
$tbl.PreferredWidth = 567.55
